$wb = $excel.ActiveWorkbook

# --- Sheet "Sheet" (matrix): update row 2 (course 1) and column B (course 1) ---
$ws1 = $wb.Worksheets.Item("Sheet")

# Row 2 updates (course 1 row)
$ws1.Range("B2").Value = 16
$ws1.Range("C2").Value = 6
$ws1.Range("E2").Value = 10
$ws1.Range("G2").Value = 7
$ws1.Range("H2").Value = 5
$ws1.Range("J2").Value = 6
$ws1.Range("K2").Value = 3
$ws1.Range("N2").Value = 3
$ws1.Range("O2").Value = 4
$ws1.Range("Q2").Value = 2
$ws1.Range("R2").Value = 5
$ws1.Range("S2").Value = 5
$ws1.Range("T2").Value = 3
$ws1.Range("U2").Value = 6

# Column B updates (course 1 column), mirroring the symmetric matrix
$ws1.Range("B3").Value = 6
$ws1.Range("B5").Value = 10
$ws1.Range("B7").Value = 7
$ws1.Range("B8").Value = 5
$ws1.Range("B10").Value = 6
$ws1.Range("B11").Value = 3
$ws1.Range("B14").Value = 3
$ws1.Range("B15").Value = 4
$ws1.Range("B17").Value = 2
$ws1.Range("B18").Value = 5
$ws1.Range("B19").Value = 5
$ws1.Range("B20").Value = 3
$ws1.Range("B21").Value = 6

# --- Sheet "SomeDataSheet": header casing fix + updated registered count ---
$ws2 = $wb.Worksheets.Item("SomeDataSheet")
$ws2.Range("B1").Value = "Number registered"
$ws2.Range("B2").Value = 16
